# Generate Report for Handback
#
# For the "4edd863f-447f-4bf2-a538-500dc46c8b47" file, a new handback was
# received (but it was generated from a stale source revision), so the
# report needs to record the handback target file, the handback xliff
# file name, the handback datetime and an error message, for both the
# zh-cn and de-de language sheets. The Error Detail column is also
# widened so the long message is readable.

$wb = $excel.ActiveWorkbook

$latestHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e8998df8ac893305c8be2752d0b1cff32b7f9078/e2e/4edd863f-447f-4bf2-a538-500dc46c8b47.md"
$targetFileDisplay = "4edd863f-447f-4bf2-a538-500dc46c8b47.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1de6a99d085a25e89fc90f291fa5de9a8910a451/e2e/4edd863f-447f-4bf2-a538-500dc46c8b47.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e8998df8ac893305c8be2752d0b1cff32b7f9078/e2e/4edd863f-447f-4bf2-a538-500dc46c8b47.md."

function Set-HandbackRow6($ws, $handbackXlf, $handbackDateTime) {
    # Latest Target File: becomes a hyperlink to the handed-back markdown file
    $ws.Hyperlinks.Add($ws.Range("I6"), $latestHandbackUrl, "", "", $targetFileDisplay)
    $ws.Range("I6").Style = "HyperLink"
    $ws.Range("I6").Font.Underline = 2
    $ws.Range("I6").Font.Color = 15570276

    # Latest Handback File
    $ws.Range("J6").Value = $handbackXlf

    # Latest Handback DateTime
    $ws.Range("K6").Value = $handbackDateTime

    # Error Detail
    $ws.Range("P6").Value = $errorDetail

    # Widen the Error Detail column so the message is readable
    $ws.Range("P1").ColumnWidth = 39.14
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow6 $wsZhCn "4edd863f-447f-4bf2-a538-500dc46c8b47.295786a8bc8971c910be1af98ac8ab588afb5e8f.zh-cn.xlf" "2016-08-25 06:42:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackRow6 $wsDeDe "4edd863f-447f-4bf2-a538-500dc46c8b47.295786a8bc8971c910be1af98ac8ab588afb5e8f.de-de.xlf" "2016-08-25 06:42:17"
